# Apply cryptocurrency price/volume updates scraped on Thu Apr 11 19:57:36 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.344.12"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "'3.512.74"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'604.45"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").Value = "'174.68"
$ws.Range("E6").Value = "  +2.24%  "
$ws.Range("E7").Value = "  -1.14%  "
$ws.Range("D8").Value = "'3.506.96"
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").Value = "'0.193"
$ws.Range("E10").Value = "  -3.70%  "
$ws.Range("D11").Value = "'7.20"
$ws.Range("E11").Value = "  +8.21%  "
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("D13").Value = "'46.15"
$ws.Range("E13").Value = "  -2.44%  "
$ws.Range("E14").Value = "  -1.53%  "
$ws.Range("D15").Value = "'4.073.18"
$ws.Range("E15").Value = "  -0.19%  "
$ws.Range("D16").Value = "'8.28"
$ws.Range("E16").Value = "  -0.85%  "
$ws.Range("D17").Value = "'608.45"
$ws.Range("E17").Value = "  -1.31%  "
$ws.Range("D18").Value = "'3.524.84"
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("D19").Value = "'70.442.77"
$ws.Range("E19").Value = "  +0.79%  "
$ws.Range("E20").Value = "  +0.84%  "
$ws.Range("D21").Value = "'17.40"
$ws.Range("E21").Value = "  +0.66%  "
$ws.Range("E22").Value = "  -0.76%  "
$ws.Range("D23").Value = "'9.00"
$ws.Range("E23").Value = "  -11.25%  "
$ws.Range("D24").Value = "'98.37"
$ws.Range("E24").Value = "  +2.41%  "
$ws.Range("D25").Value = "'15.56"
$ws.Range("E25").Value = "  -1.41%  "
$ws.Range("D26").Value = "'3.72"
$ws.Range("E26").Value = "  -3.92%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("E28").Value = "  -1.91%  "
$ws.Range("D29").Value = "'33.81"
$ws.Range("E29").Value = "  +1.83%  "
$ws.Range("D30").Value = "'8.99"
$ws.Range("E30").Value = "  -2.79%  "
$ws.Range("E31").Value = "  -3.48%  "
$ws.Range("D32").Value = "'8.03"
$ws.Range("E32").Value = "  -4.72%  "
$ws.Range("D33").Value = "'639.21"
$ws.Range("E33").Value = "  +13.65%  "
$ws.Range("E34").Value = "  -4.35%  "
$ws.Range("E35").Value = "  -1.95%  "
$ws.Range("D36").Value = "'3.57"
$ws.Range("E36").Value = "  +1.33%  "
$ws.Range("D37").Value = "'0.0993"
$ws.Range("E37").Value = "  -2.04%  "
$ws.Range("D38").Value = "'10.75"
$ws.Range("E38").Value = "  -0.19%  "
$ws.Range("D39").Value = "'0.0474"
$ws.Range("E39").Value = "  +5.21%  "
$ws.Range("D40").Value = "'56.72"
$ws.Range("E40").Value = "  -0.65%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.16%  "
$ws.Range("E42").Value = "  +0.81%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "'3.373.20"
$ws.Range("E43").Value = "  +1.20%  "
$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D44").Value = "'0.0₃0740"
$ws.Range("E44").Value = "  +5.02%  "
$ws.Range("E45").Value = "  -5.73%  "
$ws.Range("D46").Value = "'32.14"
$ws.Range("E46").Value = "  -2.73%  "
$ws.Range("E47").Value = "  -0.76%  "
$ws.Range("D48").Value = "'2.55"
$ws.Range("E48").Value = "  -2.70%  "
$ws.Range("E49").Value = "  +0.54%  "
$ws.Range("D50").Value = "'132.63"
$ws.Range("E50").Value = "  -2.80%  "
$ws.Range("E51").Value = "  -0.01%  "
